$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.183.30"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.989.19"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.12"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.11"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.40"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("E12").Value = "  -4.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.35"
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.463.63"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.74"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.977.50"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.03"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.212.20"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +5.85%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.68"
$ws.Range("E21").Value = "  -5.47%  "
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.54"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.28"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.181"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.86"
$ws.Range("E27").Value = "  -3.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.55"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.109"
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.38"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.76"
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("E34").Value = "  +10.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.94"
$ws.Range("E35").Value = "  -3.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0448"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.21"
$ws.Range("E38").Value = "  -6.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.92"
$ws.Range("E39").Value = "  -5.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.98"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.70"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.83"
$ws.Range("E43").Value = "  -4.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.03"
$ws.Range("E44").Value = "  +7.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.126.09"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.40"
$ws.Range("E47").Value = "  -4.40%  "
$ws.Range("E48").Value = "  -5.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.249"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("E51").Value = "  -1.91%  "
